$cmds = Get-Command
Write-Host $cmds
